# Update cryptos list: refresh Price (D) and Volume(1h) (E) values
# (leading apostrophe forces plain-numeric-looking D values to remain text,
#  matching the original inlineStr cell type instead of being auto-converted
#  to a number by Excel)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.335.21"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.864.95"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'239.43"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'42.88"
$ws.Range("E8").Value = "  +7.57%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "2.133.43"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "'11.56"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "1.841.19"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "35.338.93"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "'70.25"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'241.47"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'12.28"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'4.75"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "'169.59"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'1.91"
$ws.Range("E26").Value = "  +25.86%  "
$ws.Range("D27").Value = "'8.21"
$ws.Range("E27").Value = "  +5.48%  "
$ws.Range("D28").Value = "'17.78"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "'1.84"
$ws.Range("E33").Value = "  +27.28%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("D36").Value = "'0.816"
$ws.Range("E36").Value = "  +17.52%  "
$ws.Range("E37").Value = "  +6.29%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("D40").Value = "'91.15"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "1.348.54"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'15.22"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  +15.33%  "
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "'12.99"
$ws.Range("E45").Value = "  +54.43%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  +5.68%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "2.051.91"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").Value = "'0.0686"
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("E51").Value = "  -0.78%  "